$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 4 (phone +5511920075911 / DDD 11 / 2024-10-25); subsequent rows shift up.
$ws.Rows.Item(4).Delete()
